$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 2: URL  -> http://ibm.com/...        => http://linuxforhealth.org/...
$ws.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/ms-drg-of-payment"

# Row 3: Version  7.0.0 => 8.0.0
$ws.Range("B3").Value = "8.0.0"

# Row 8: Date  2022-09-08T16:11:15+00:00 => 2022-11-10T16:00:46+00:00
$ws.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Row 9: Publisher  Alvearie Team => LinuxForHealth Team
$ws.Range("B9").Value = "LinuxForHealth Team"
